$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Row 19 : Armitage - "Vernal behaviour of the yellow-bellied marmot"
# ---------------------------------------------------------------------------
CopyFormat "A43" "A19"
CopyFormat "I9"  "B19"
CopyFormat "D10" "C19"
CopyFormat "I9"  "D19"
CopyFormat "I9"  "E19"
CopyFormat "F2"  "F19"
CopyFormat "I9"  "G19"
CopyFormat "I9"  "H19"
CopyFormat "I9"  "I19"

$ws.Range("A19").Value = "Vernal behaviour of the yellow-bellied marmot"
$ws.Range("B19").Value = "Anim. Behav."
$ws.Range("C19").Value = 1965
$ws.Range("D19").Value = "Armitage"
$ws.Range("E19").Value = "YELLOWSTONE POPULATION"
$ws.Range("F19").Value = "yes"
$ws.Range("G19").Value = "Reproductive and agonistic behaviour in the first 5 weeks post-hibernation"
$ws.Range("H19").Value = "behavioral observation in natural conditions"
$ws.Range("I19").Value = "yes"
$ws.Range("J19").Value = "yes"
$ws.Range("K19").Value = "yes"

# ---------------------------------------------------------------------------
# Row 20 : Blumstein, Nicodemus and Zugmeyer (2004) - J. Mamm.
# ---------------------------------------------------------------------------
CopyFormat "F7" "F20"
CopyFormat "B2" "D20"

$ws.Hyperlinks.Add($ws.Range("A20"), "https://sites.lifesci.ucla.edu/eeb-rmbl-marmots/wp-content/uploads/sites/190/2020/03/Blumstein_etal_2004_JM.pdf", "", "", "https://sites.lifesci.ucla.edu/eeb-rmbl-marmots/wp-content/uploads/sites/190/2020/03/Blumstein_etal_2004_JM.pdf") | Out-Null
$ws.Range("A20").Value = "Yellow-bellied marmots (Marmota flaviventris) hibernate socially"
$ws.Range("A20").Style = "Hyperlink"

$ws.Range("B20").Value = "J. Mamm."
$ws.Range("C20").Value = 2004
$ws.Range("D20").Value = "Blumstein, Nicodemus and Zugmeyer"
$ws.Range("F20").Value = "no"
$ws.Range("I20").Value = "yes"
$ws.Range("J20").Value = "yes"
$ws.Range("K20").Value = "yes"

# ---------------------------------------------------------------------------
# Row 21 : Blumstein (2009) - J. Mamm.
# ---------------------------------------------------------------------------
CopyFormat "F7" "F21"

$ws.Hyperlinks.Add($ws.Range("A21"), "https://sites.lifesci.ucla.edu/eeb-rmbl-marmots/wp-content/uploads/sites/190/2020/03/Blumstein_2009_JM.pdf", "", "", "https://sites.lifesci.ucla.edu/eeb-rmbl-marmots/wp-content/uploads/sites/190/2020/03/Blumstein_2009_JM.pdf") | Out-Null
$ws.Range("A21").Value = "Social effects on emergence from hibernation in yellow-bellied marmots"
$ws.Range("A21").Style = "Hyperlink"

$ws.Range("B21").Value = "J. Mamm."
$ws.Range("C21").Value = 2009
$ws.Range("D21").Value = "Blumstein"
$ws.Range("F21").Value = "no"
$ws.Range("I21").Value = "yes"
$ws.Range("J21").Value = "yes"
$ws.Range("K21").Value = "yes"

# ---------------------------------------------------------------------------
# Row 22 : Monclus, Pang, Blumstein (2014) - Evolutionary Ecology
# ---------------------------------------------------------------------------
CopyFormat "F7" "F22"
CopyFormat "B2" "B22"
CopyFormat "B2" "D22"
CopyFormat "J20" "I22"
CopyFormat "J20" "J22"
CopyFormat "J20" "K22"

$ws.Hyperlinks.Add($ws.Range("A22"), "https://sites.lifesci.ucla.edu/eeb-rmbl-marmots/wp-content/uploads/sites/190/2020/03/Monclus_etal_2014_EvolEcol.pdf", "", "", "https://sites.lifesci.ucla.edu/eeb-rmbl-marmots/wp-content/uploads/sites/190/2020/03/Monclus_etal_2014_EvolEcol.pdf") | Out-Null
$ws.Range("A22").Value = " Yellow-bellied marmots do not compensate for a late start: the role of maternal investment in shaping life-history trajectories"
$ws.Range("A22").Style = "Hyperlink"

$ws.Range("B22").Value = "Evolutionary Ecology"
$ws.Range("C22").Value = 2014
$ws.Range("D22").Value = "Monclus, Pang, Blumstein"
$ws.Range("F22").Value = "no"
$ws.Range("I22").Value = "yes"
$ws.Range("J22").Value = "yes"
$ws.Range("K22").Value = "yes"

# ---------------------------------------------------------------------------
# Housekeeping: Excel's save-time style-table GC re-pointed these unchanged
# cells at equivalent (deduplicated) cellXfs entries. Mirror the same visual
# style using an existing equivalent entry already on the sheet.
# ---------------------------------------------------------------------------
CopyFormat "A3" "M2"
$ws.Range("M2").Value = "Point of the paper"

CopyFormat "I9" "K18"
$ws.Range("K18").Value = "yes"

$wb.Application.CalculateFull()
